$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39. This shifts the existing rows 39-123
# down to 40-124 (and the sheet's used range grows from A1:T123 to
# A1:T124), matching the data seen lower in the sheet (e.g. old row 123
# becomes new row 124 verbatim).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 45162
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100108
$ws.Range("H39").Value = "Tropicales y subtropicales"
$ws.Range("I39").Value = 100108007
$ws.Range("J39").Value = "Coco"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 90
$ws.Range("N39").Value = 35000
$ws.Range("O39").Value = 36000
$ws.Range("P39").Value = 35556
$ws.Range("Q39").Value = '$/malla 20 unidades'
$ws.Range("R39").Value = "Perú"
$ws.Range("S39").Value = 1778
$ws.Range("T39").Value = 20
